$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Range("D2")
$cell.Value = "'246.89"
$cell.Style = "Normal"

# Row 3
$cell = $ws.Range("D3")
$cell.Value = "'26.29"
$cell.Style = "Normal"

# Row 4
$cell = $ws.Range("D4")
$cell.Value = "'5.084"
$cell.Style = "Normal"

# Row 5
$cell = $ws.Range("D5")
$cell.Value = "'0.05620"
$cell.Style = "Normal"

# Row 6
$cell = $ws.Range("D6")
$cell.Value = "'6.492"
$cell.Style = "Normal"

# Row 7
$ws.Range("B7").Value = 'MXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$cell = $ws.Range("D7")
$cell.Value = "'0.8129"
$cell.Style = "Normal"
$ws.Range("E7").Value = '6MXTokenMX'

# Row 8
$ws.Range("B8").Value = 'FTXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$cell = $ws.Range("D8")
$cell.Value = "'0.8470"
$cell.Style = "Normal"
$ws.Range("E8").Value = '7FTXTokenFTT'

# Row 9
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$cell = $ws.Range("D9")
$cell.Value = "'0.03164"
$cell.Style = "Normal"
$ws.Range("E9").Value = '8LiechtensteinCryptoassetsExchangeLCX'

# Row 10
$ws.Range("B10").Value = 'BitrueCoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$cell = $ws.Range("D10")
$cell.Value = "'0.02821"
$cell.Style = "Normal"
$ws.Range("E10").Value = '9BitrueCoinBTR'

# Row 11
$ws.Range("B11").Value = 'BitMartToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$cell = $ws.Range("D11")
$cell.Value = "'0.09405"
$cell.Style = "Normal"
$ws.Range("E11").Value = '10BitMartTokenBMX'

# Row 12
$ws.Range("B12").Value = 'BitForexToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$cell = $ws.Range("D12")
$cell.Value = "'0.001521"
$cell.Style = "Normal"
$ws.Range("E12").Value = '11BitForexTokenBF'

# Row 13
$ws.Range("B13").Value = 'One'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$cell = $ws.Range("D13")
$cell.Value = "'0.0005982"
$cell.Style = "Normal"
$ws.Range("E13").Value = '12OneONE'

# Row 14
$ws.Range("B14").Value = 'TigerCash'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$cell = $ws.Range("D14")
$cell.Value = "'0.006239"
$cell.Style = "Normal"
$ws.Range("E14").Value = '13TigerCashTCH'

# Row 15
$ws.Range("B15").Value = 'LEO'
$ws.Range("C15").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$cell = $ws.Range("D15")
$cell.Value = "'3.573"
$cell.Style = "Normal"
$ws.Range("E15").Value = '14LEOLEO'

# Row 16
$ws.Range("B16").Value = 'GateToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$cell = $ws.Range("D16")
$cell.Value = "'3.057"
$cell.Style = "Normal"
$ws.Range("E16").Value = '15GateTokenGT'

# Row 17
$ws.Range("B17").Value = 'BTSEToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$cell = $ws.Range("D17")
$cell.Value = "'2.118"
$cell.Style = "Normal"
$ws.Range("E17").Value = '16BTSETokenBTSE'

# Row 18
$ws.Range("B18").Value = 'BitpandaEcosystemToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$cell = $ws.Range("D18")
$cell.Value = "'0.3181"
$cell.Style = "Normal"
$ws.Range("E18").Value = '17BitpandaEcosystemTokenBEST'

# Row 19
$ws.Range("B19").Value = 'WazirX'
$ws.Range("C19").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$cell = $ws.Range("D19")
$cell.Value = "'0.1347"
$cell.Style = "Normal"
$ws.Range("E19").Value = '18WazirXWRX'

# Row 20
$cell = $ws.Range("D20")
$cell.Value = "'0.06957"
$cell.Style = "Normal"

# Row 22
$cell = $ws.Range("D22")
$cell.Value = "'3.758"
$cell.Style = "Normal"

# Row 23
$cell = $ws.Range("D23")
$cell.Value = "'0.04680"
$cell.Style = "Normal"

# Row 25
$cell = $ws.Range("D25")
$cell.Value = "'0.001252"
$cell.Style = "Normal"

# Row 26
$cell = $ws.Range("D26")
$cell.Value = "'0.004622"
$cell.Style = "Normal"

# Row 27
$cell = $ws.Range("D27")
$cell.Value = "'0.00009603"
$cell.Style = "Normal"

# Row 28
$cell = $ws.Range("D28")
$cell.Value = "'0.0001939"
$cell.Style = "Normal"

# Row 40
$cell = $ws.Range("D40")
$cell.Value = "'0.03664"
$cell.Style = "Normal"

# Row 41
$cell = $ws.Range("D41")
$cell.Value = "'0.006112"
$cell.Style = "Normal"
$ws.Range("E41").Value = '40KickTokenKICKBestin24h'

# Row 42
$cell = $ws.Range("D42")
$cell.Value = "'0.1059"
$cell.Style = "Normal"

# Row 43
$cell = $ws.Range("D43")
$cell.Value = "'0.002598"
$cell.Style = "Normal"

# Row 44
$cell = $ws.Range("D44")
$cell.Value = "'0.008686"
$cell.Style = "Normal"
$ws.Range("E44").Value = '43LocalTradersLCT'

# Row 45
$cell = $ws.Range("D45")
$cell.Value = "'0.00005295"
$cell.Style = "Normal"

# Row 47
$ws.Range("E47").Value = '46CoinbaseStockTokenCOINWorstin24h'

# Row 48
$cell = $ws.Range("D48")
$cell.Value = "'0.002066"
$cell.Style = "Normal"

# Row 49
$cell = $ws.Range("D49")
$cell.Value = "'0.00002101"
$cell.Style = "Normal"

# Row 50
$cell = $ws.Range("D50")
$cell.Value = "'0.0002001"
$cell.Style = "Normal"
